# Rename wc_lang DfbaNetReaction -> DfbaObjReaction, DfbaNetSpecies -> DfbaObjSpecies
# Reflected in the workbook as:
#   - sheet "dFBA net reactions" -> "dFBA objective reactions"
#   - sheet "dFBA net species"   -> "dFBA objective species"
#   - "Parameters" sheet column header "dFBA net reaction" -> "dFBA objective reaction"

$wb = $excel.ActiveWorkbook

# --- Rename the two sheets ---
$netReactions = $wb.Worksheets.Item("dFBA net reactions")
$netReactions.Name = "dFBA objective reactions"

$netSpecies = $wb.Worksheets.Item("dFBA net species")
$netSpecies.Name = "dFBA objective species"

# --- Update the column header text on the Parameters sheet ---
$parameters = $wb.Worksheets.Item("Parameters")
$parameters.Range("C1").Value = "dFBA objective reaction"

# --- View/selection bookkeeping: move the active tab onto the renamed
#     "dFBA objective species" sheet, and leave the Parameters sheet's
#     cursor on E6 (was A1) ---
$parameters.Range("E6").Select()
$netSpecies.Activate()
